# Cleaned up code and naming convention.
# Edited default_animal_config.xlsx to change all "n/a" placeholder values
# to 0. Code now works for all animal types.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 23 = fattening_pig, Row 24 = breeding_pig
# Row 25 = laying_hen, Row 26 = pullet, Row 27 = broiler_chicken
$rowsSolid  = @(23, 24)
$rowsLiquid = @(25, 26, 27)

foreach ($r in $rowsSolid) {
    # methane_potential_solid / stdev_methane_potential_solid
    $ws.Range("M$r").Value = 0
    $ws.Range("N$r").Value = 0
    # Unit Straw
    $ws.Range("R$r").Value = "L / kg VS"
    # dry_weight_solid / stdev_dry_weight_solid / organic_dry_weight_solid / stdev_organic_dry_weight_solid
    $ws.Range("Y$r").Value = 0
    $ws.Range("Z$r").Value = 0
    $ws.Range("AA$r").Value = 0
    $ws.Range("AB$r").Value = 0
}

foreach ($r in $rowsLiquid) {
    # methane_potential_liquid / stdev_methane_potential_liquid
    $ws.Range("K$r").Value = 0
    $ws.Range("L$r").Value = 0
    # Unit Straw
    $ws.Range("R$r").Value = "L / kg VS"
    # dry_weight_liquid / stdev_dry_weight_liquid / organic_dry_weight_liquid / stdev_organic_dry_weight_liquid
    $ws.Range("U$r").Value = 0
    $ws.Range("V$r").Value = 0
    $ws.Range("W$r").Value = 0
    $ws.Range("X$r").Value = 0
}

# Match the final selection left behind in the authored workbook.
$ws.Range("S31").Select()
